$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    [PSCustomObject]@{ Row = 2; D = "26.934.48"; E = "  -0.59%  " }
    [PSCustomObject]@{ Row = 3; D = "1.863.63"; E = "  -0.23%  " }
    [PSCustomObject]@{ Row = 4; D = "0.9990"; E = "  -0.07%  " }
    [PSCustomObject]@{ Row = 5; D = "304.85"; E = "  -0.72%  " }
    [PSCustomObject]@{ Row = 6; D = "0.9992"; E = "  -0.08%  " }
    [PSCustomObject]@{ Row = 7; D = "0.5074"; E = "  -0.28%  " }
    [PSCustomObject]@{ Row = 8; D = "0.3641"; E = "  -2.68%  " }
    [PSCustomObject]@{ Row = 9; D = "0.07181"; E = "  +0.66%  " }
    [PSCustomObject]@{ Row = 10; D = "0.8959"; E = "  +1.24%  " }
    [PSCustomObject]@{ Row = 11; D = "20.82"; E = "  +1.03%  " }
    [PSCustomObject]@{ Row = 12; D = "1.870.39"; E = "  +0.03%  " }
    [PSCustomObject]@{ Row = 13; D = "0.07494"; E = "  -0.40%  " }
    [PSCustomObject]@{ Row = 14; D = "92.55"; E = "  +3.92%  " }
    [PSCustomObject]@{ Row = 15; D = "5.235"; E = "  -1.40%  " }
    [PSCustomObject]@{ Row = 16; D = "0.9994"; E = "  -0.06%  " }
    [PSCustomObject]@{ Row = 17; D = "0.000008498"; E = "  +0.41%  " }
    [PSCustomObject]@{ Row = 18; D = "14.20"; E = "  +0.69%  " }
    [PSCustomObject]@{ Row = 19; D = "0.9997"; E = "  -0.03%  " }
    [PSCustomObject]@{ Row = 20; D = "26.970.76"; E = "  -0.66%  " }
    [PSCustomObject]@{ Row = 21; D = "5.040"; E = "  -0.08%  " }
    [PSCustomObject]@{ Row = 22; D = "2.100.44"; E = "  -0.62%  " }
    [PSCustomObject]@{ Row = 23; D = "10.38"; E = "  -1.59%  " }
    [PSCustomObject]@{ Row = 24; D = "6.413"; E = "  -0.90%  " }
    [PSCustomObject]@{ Row = 25; D = "147.46"; E = "  -1.44%  " }
    [PSCustomObject]@{ Row = 26; D = "1.795"; E = "  -3.02%  " }
    [PSCustomObject]@{ Row = 27; D = "17.88"; E = "  -0.25%  " }
    [PSCustomObject]@{ Row = 28; D = "2.066"; E = "  -1.49%  " }
    [PSCustomObject]@{ Row = 29; D = "113.14"; E = "  +0.36%  " }
    [PSCustomObject]@{ Row = 30; D = "4.684"; E = $null }
    [PSCustomObject]@{ Row = 31; D = "4.686"; E = "  -0.01%  " }
    [PSCustomObject]@{ Row = 32; D = "0.09262"; E = "  +2.83%  " }
    [PSCustomObject]@{ Row = 33; D = "0.05103"; E = "  -0.43%  " }
    [PSCustomObject]@{ Row = 34; D = "0.7514"; E = "  +2.53%  " }
    [PSCustomObject]@{ Row = 35; D = "2.991"; E = "  -3.21%  " }
    [PSCustomObject]@{ Row = 36; D = "1.152"; E = "  -0.60%  " }
    [PSCustomObject]@{ Row = 37; D = $null; E = "  +7.52%  " }
    [PSCustomObject]@{ Row = 38; D = "2.548"; E = "  +1.79%  " }
    [PSCustomObject]@{ Row = 39; D = "0.02003"; E = "  -2.10%  " }
    [PSCustomObject]@{ Row = 40; D = "0.5539"; E = "  +4.20%  " }
    [PSCustomObject]@{ Row = 41; D = "1.073"; E = "  -0.49%  " }
    [PSCustomObject]@{ Row = 42; D = "118.12"; E = "  +0.85%  " }
    [PSCustomObject]@{ Row = 43; D = "6.534"; E = "  -0.29%  " }
    [PSCustomObject]@{ Row = 44; D = "8.498"; E = "  +2.29%  " }
    [PSCustomObject]@{ Row = 45; D = "0.1471"; E = "  +0.19%  " }
    [PSCustomObject]@{ Row = 46; D = "0.4687"; E = $null }
    [PSCustomObject]@{ Row = 47; D = "0.9989"; E = "  -0.08%  " }
    [PSCustomObject]@{ Row = 48; D = $null; E = "  +0.51%  " }
    [PSCustomObject]@{ Row = 49; D = "1.567"; E = "  +0.34%  " }
    [PSCustomObject]@{ Row = 50; D = "36.83"; E = "  +1.01%  " }
    [PSCustomObject]@{ Row = 51; D = "63.15"; E = "  -1.93%  " }
)


foreach ($r in $rows) {
    $rowNum = $r.Row

    if ($null -ne $r.D) {
        $dCell = $ws.Cells.Item($rowNum, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $r.D
        $dCell.Style = "Normal"
    }

    if ($null -ne $r.E) {
        $eCell = $ws.Cells.Item($rowNum, 5)
        $eCell.NumberFormat = "@"
        $eCell.Value = $r.E
        $eCell.Style = "Normal"
    }
}
